$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new blank row above row 3, which pushes the existing
# "Who said 'Give me liberty...'" question row down from row 3 to row 4.
$ws.Rows.Item(3).Insert()
# Fully clear the newly-inserted blank row 3's leftover cell shells so it
# does not linger as an empty row in the sheet (only rows 4.. are real data).
$ws.Range("A3:K3").Clear()

# A far-away, untouched cell used purely as a "General format" stamp so we
# can paste its (default) number format onto a target cell before writing a
# literal number into it. Columns B:K carry a Text ("@") number format by
# default (inherited from the sheet's <col> styles), and this engine stores
# any value typed into a Text-formatted cell as a string - so numeric cells
# need their format reset to General first, exactly like typing into a
# fresh/general cell in Excel.
$blankFormat = $ws.Cells.Item(500, 500)
$blankFormat.Copy()

# --- Append 6 more "Civil War" quiz rows (rows 5-10), each repeating the
# same answer-option / year / flag pattern as the existing rows.
$questions = @(
    "When was the Civil War 3?",
    "When was the Civil War 4?",
    "When was the Civil War 5?",
    "When was the Civil War 6?",
    "When was the Civil War 7?",
    "When was the Civil War 8?"
)

$row = 5
foreach ($q in $questions) {
    $ws.Cells.Item($row, 1).Value = $q

    $ws.Cells.Item($row, 2).Value = "R"
    $ws.Cells.Item($row, 2).HorizontalAlignment = -4108  # xlCenter

    $ws.Cells.Item($row, 3).Value = "E"
    $ws.Cells.Item($row, 3).HorizontalAlignment = -4108  # xlCenter

    $numCells = @(
        @{ Col = 4;  Val = 1941 },
        @{ Col = 5;  Val = 0 },
        @{ Col = 6;  Val = 1942 },
        @{ Col = 7;  Val = 1 },
        @{ Col = 8;  Val = 1943 },
        @{ Col = 9;  Val = 0 },
        @{ Col = 10; Val = 1944 },
        @{ Col = 11; Val = 0 }
    )
    foreach ($nc in $numCells) {
        $cell = $ws.Cells.Item($row, $nc.Col)
        $cell.PasteSpecial(-4122)   # xlPasteFormats: reset to General (no text coercion)
        $cell.Value = $nc.Val
        $cell.NumberFormat = "@"   # restore the Text display format used by the rest of the table
        $cell.HorizontalAlignment = -4108  # xlCenter
    }

    $row = $row + 1
}

# --- Row 11 stays empty; row 12 gets a single stray cell (looks like a
# leftover "`" marker typed while reviewing results).
$ws.Cells.Item(12, 1).Value = "``"

# --- Restore the active selection to match where the editor ended up.
$ws.Range("B10").Select()
